$wb = $excel.ActiveWorkbook

# --- Sheet 1: "VENTAS POR GRUPO" (columns A:R) ---
$ws1 = $wb.Worksheets.Item(1)

# Insert a new row at position 6, shifting the existing rows (and their
# formatting) down by one.
$ws1.Rows.Item(6).Insert()

$ws1.Cells.Item(6, 1).Value = "VACA PANCHI CAROLINA"
$ws1.Cells.Item(6, 2).Value = "GRANJA VANEGAS MARCELA"
for ($col = 3; $col -le 18; $col++) {
    $ws1.Cells.Item(6, $col).Value = 0
}

# Update the totals row (now row 12) label from "0 de 9" to "0 de 10"
for ($col = 3; $col -le 18; $col++) {
    $ws1.Cells.Item(12, $col).Value = "0 de 10"
}

# --- Sheet 2: "VENTA MENSUAL" (columns A:G) ---
$ws2 = $wb.Worksheets.Item(2)

# Insert a new row at position 6, shifting the existing rows (and their
# formatting) down by one.
$ws2.Rows.Item(6).Insert()

$ws2.Cells.Item(6, 1).Value = "VACA PANCHI CAROLINA"
$ws2.Cells.Item(6, 2).Value = "GRANJA VANEGAS MARCELA"
for ($col = 3; $col -le 7; $col++) {
    $ws2.Cells.Item(6, $col).Value = 0
}
